$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1), matching the formatting of the other
# header cells (copy format from G1, then set the text).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data values for the "Save" column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
